$d = $word.ActiveDocument

# The WordprocessingML fragment for a new, empty "ListParagraph" styled
# paragraph indented 1080 twips (54pt) - no numbering, no bookmark, no run.
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/></w:pPr></w:p>'

function Find-BookmarkParaIndex($doc) {
    # The target is the numbered-list paragraph carrying the "_GoBack"
    # bookmark - i.e. the last paragraph in the document with an active
    # list/numbering format.
    $n = $doc.Paragraphs.Count
    for ($i = $n; $i -ge 1; $i--) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.ListFormat.ListType -ne 0) {
            return $i
        }
    }
    throw "Could not locate the target numbered list paragraph"
}

# --- Insert two new empty ListParagraph paragraphs immediately BEFORE the
#     bookmark paragraph. Re-resolve the paragraph index/position fresh
#     each time since prior COM Range/Paragraph handles do not reseat
#     themselves after a structural edit. ---
for ($k = 0; $k -lt 2; $k++) {
    $idx = Find-BookmarkParaIndex $d
    $insPos = $d.Paragraphs($idx).Range.Start
    $r = $d.Range($insPos, $insPos)
    [void]$r.InsertXML($xmlFrag)
}

# --- Insert one new empty ListParagraph paragraph immediately AFTER the
#     bookmark paragraph. InsertParagraphAfter() creates a correctly
#     positioned trailing paragraph; InsertXML over that fresh paragraph's
#     own range then replaces it with the exact target markup (no
#     inherited numbering, no stray run). ---
$idx = Find-BookmarkParaIndex $d
$target = $d.Paragraphs($idx)
[void]$target.Range.InsertParagraphAfter()
$np = $d.Paragraphs($idx + 1)
[void]$np.Range.InsertXML($xmlFrag)
